$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (B and C columns)
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 1

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 2

$ws.Range("B4").Value = 26
$ws.Range("C4").Value = 1

# Add new row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 26
$ws.Range("C5").Value = 2

# Update selection to match target (active cell A2, selected range A2:C4)
$ws.Range("A2:C4").Select()
